$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Risky" Price cells are single-dot decimals (e.g. "579.13") that Excel would
# otherwise auto-convert to a floating point number (losing the exact source
# text / introducing binary rounding). Mark each one as Text *before* assigning
# its value so it is stored verbatim, exactly like the other inline-string cells.
$riskyPriceCells = @('D4', 'D5', 'D6', 'D7', 'D8', 'D9', 'D11', 'D14', 'D16', 'D18', 'D19', 'D20', 'D21', 'D22', 'D23', 'D24', 'D26', 'D27', 'D29', 'D30', 'D31', 'D32', 'D33', 'D34', 'D35', 'D36', 'D37', 'D38', 'D39', 'D40', 'D41', 'D42', 'D43', 'D44', 'D46', 'D47', 'D48', 'D49', 'D50', 'D51')
foreach ($addr in $riskyPriceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '67.972.57'
$ws.Range("E2").Value = '  +0.67%  '

$ws.Range("D3").Value = '3.245.51'
$ws.Range("E3").Value = '  -0.27%  '

$ws.Range("D4").Value = '0.998'
$ws.Range("E4").Value = '  -0.19%  '

$ws.Range("D5").Value = '579.13'
$ws.Range("E5").Value = '  -0.01%  '

$ws.Range("D6").Value = '182.55'
$ws.Range("E6").Value = '  +0.66%  '

$ws.Range("D7").Value = '0.998'
$ws.Range("E7").Value = '  -0.19%  '

$ws.Range("D8").Value = '0.596'
$ws.Range("E8").Value = '  -0.39%  '

$ws.Range("D9").Value = '0.134'
$ws.Range("E9").Value = '  +1.00%  '

$ws.Range("E10").Value = '  -1.65%  '

$ws.Range("D11").Value = '0.418'
$ws.Range("E11").Value = '  +0.65%  '

$ws.Range("D12").Value = '3.791.14'
$ws.Range("E12").Value = '  -0.75%  '

$ws.Range("E13").Value = '  -0.38%  '

$ws.Range("D14").Value = '28.16'
$ws.Range("E14").Value = '  -0.65%  '

$ws.Range("D15").Value = '67.869.40'
$ws.Range("E15").Value = '  +0.52%  '

$ws.Range("D16").Value = '0.0000171'
$ws.Range("E16").Value = '  +1.69%  '

$ws.Range("D17").Value = '3.201.93'
$ws.Range("E17").Value = '  -1.68%  '

$ws.Range("D18").Value = '5.81'
$ws.Range("E18").Value = '  -0.37%  '

$ws.Range("D19").Value = '13.49'
$ws.Range("E19").Value = '  -0.15%  '

$ws.Range("D20").Value = '391.26'
$ws.Range("E20").Value = '  +4.41%  '

$ws.Range("D21").Value = '7.67'
$ws.Range("E21").Value = '  +0.60%  '

$ws.Range("B22").Value = 'Litecoin'
$ws.Range("C22").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D22").Value = '71.62'
$ws.Range("E22").Value = '  +0.79%  '

$ws.Range("B23").Value = 'Dai'
$ws.Range("C23").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D23").Value = '0.999'
$ws.Range("E23").Value = '  -0.16%  '

$ws.Range("D24").Value = '0.516'
$ws.Range("E24").Value = '  +1.05%  '

$ws.Range("E25").Value = '  -0.03%  '

$ws.Range("D26").Value = '0.188'
$ws.Range("E26").Value = '  +3.91%  '

$ws.Range("D27").Value = '9.57'
$ws.Range("E27").Value = '  +0.05%  '

$ws.Range("E28").Value = '  +0.15%  '

$ws.Range("D29").Value = '1.97'
$ws.Range("E29").Value = '  -0.52%  '

$ws.Range("D30").Value = '5.66'
$ws.Range("E30").Value = '  +0.08%  '

$ws.Range("D31").Value = '23.01'
$ws.Range("E31").Value = '  +1.37%  '

$ws.Range("D32").Value = '7.11'
$ws.Range("E32").Value = '  +2.99%  '

$ws.Range("D33").Value = '0.999'
$ws.Range("E33").Value = '  +0.11%  '

$ws.Range("D34").Value = '1.28'
$ws.Range("E34").Value = '  +0.83%  '

$ws.Range("D35").Value = '164.49'
$ws.Range("E35").Value = '  +0.91%  '

$ws.Range("D36").Value = '1.49'
$ws.Range("E36").Value = '  -1.15%  '

$ws.Range("D37").Value = '1.91'
$ws.Range("E37").Value = '  +3.09%  '

$ws.Range("D38").Value = '0.821'
$ws.Range("E38").Value = '  -3.10%  '

$ws.Range("D39").Value = '26.56'
$ws.Range("E39").Value = '  -0.58%  '

$ws.Range("D40").Value = '4.61'
$ws.Range("E40").Value = '  -0.56%  '

$ws.Range("D41").Value = '6.51'
$ws.Range("E41").Value = '  -4.49%  '

$ws.Range("D42").Value = '2.51'
$ws.Range("E42").Value = '  -3.42%  '

$ws.Range("D43").Value = '41.36'
$ws.Range("E43").Value = '  +1.45%  '

$ws.Range("D44").Value = '0.0684'
$ws.Range("E44").Value = '  +0.75%  '

$ws.Range("D45").Value = '2.615.79'
$ws.Range("E45").Value = '  -2.99%  '

$ws.Range("D46").Value = '340.62'
$ws.Range("E46").Value = '  -3.10%  '

$ws.Range("D47").Value = '24.77'
$ws.Range("E47").Value = '  -2.22%  '

$ws.Range("D48").Value = '0.0280'
$ws.Range("E48").Value = '  -0.07%  '

$ws.Range("D49").Value = '6.32'
$ws.Range("E49").Value = '  +2.89%  '

$ws.Range("B50").Value = 'Arweave'
$ws.Range("C50").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range("D50").Value = '31.51'
$ws.Range("E50").Value = '  +0.68%  '

$ws.Range("B51").Value = 'Stellar'
$ws.Range("C51").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D51").Value = '0.102'
$ws.Range("E51").Value = '  -0.18%  '

# Drop back to the default (General) cell style now that the text is safely
# stored, so no stray style index is left on these cells.
foreach ($addr in $riskyPriceCells) {
    $ws.Range($addr).Style = "Normal"
}

